# Add three new bluetooth "customer" login rows (ids 2, 3, 4) below the
# existing seed row, mirroring how the Python login app appends a record
# each time a new device/user pairs over bluetooth.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: customer "Eyad" --------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "customer"
$ws.Range("C3").Value = "Eyad"
$ws.Range("D3").Value = "None"
$ws.Range("E3").Value = "78:46:D4:55:D8:12"
$ws.Range("F3").Value = "C:\Users\Lenovo\Downloads\download.jpg"

# --- Row 4: customer "lala" ---------------------------------------------
# (mac address entered before the display name for this record)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "customer"
$ws.Range("E4").Value = "5C-61-99-44-A2-E4"
$ws.Range("C4").Value = "lala"
$ws.Range("D4").Value = "None"
$ws.Range("F4").Value = "C:\Users\Lenovo\Downloads\download.jpg"

# --- Row 5: customer "lalal" --------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "customer"
$ws.Range("C5").Value = "lalal"
$ws.Range("D5").Value = "None"
$ws.Range("E5").Value = "5C-61-99-44-A2-E3"
$ws.Range("F5").Value = "C:\Users\Lenovo\Downloads\download.jpg"

# Column E (mac_address) re-fit itself to the new data in the source file.
$ws.Range("E1:E5").Columns.AutoFit() | Out-Null

# The saved workbook's active selection ended up on the last filled cell.
$ws.Range("F5").Select() | Out-Null
